# "Generate Report for Handback"
#
# The localization-status report gets a handback pass:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview + per-language sheets).
#   - Each per-language sheet (zh-cn, de-de) gains two new tracking columns that
#     were already reserved in the table header ("Latest Target File" / F, and
#     "Latest Handback File" / G): they get populated with the same file
#     references as the existing "Source File Name" (A) / "Latest Handoff File"
#     (D) columns, including matching hyperlinks.
#   - "Latest Handback DateTime" (H) moves from the zero-date sentinel to a
#     real timestamp, per language.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 1. Status text: Overview (B/C for rows 2-3) + each language sheet's Status
#    column (C for rows 2-3).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value2 = $newStatus
$wsOverview.Range("C2").Value2 = $newStatus
$wsOverview.Range("B3").Value2 = $newStatus
$wsOverview.Range("C3").Value2 = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("C3").Value2 = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("C3").Value2 = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: Latest Target File (F) / Latest Handback File (G), and the
#    Latest Handback DateTime (H).
# ---------------------------------------------------------------------------
$zhUrlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/6d562aac7f97325115fbb32b09ff899fb3811720/e2e/61df56d7-e540-4367-8ba3-57a9d998e113.md"
$zhUrlD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7709e2988e6f7db0ae3a6205ee195e8fbe5cfae3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/61df56d7-e540-4367-8ba3-57a9d998e113.2c65cb27da0e9d1c9c46324e0027fa7f7b8ec244.zh-cn.xlf"
$zhUrlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/6d562aac7f97325115fbb32b09ff899fb3811720/e2e/b3c9e954-dfb6-437b-b587-b10869f38c87.md"
$zhUrlD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7709e2988e6f7db0ae3a6205ee195e8fbe5cfae3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b3c9e954-dfb6-437b-b587-b10869f38c87.264bcc2d9140a76805ec1f79e3eaf0d34992c5a9.zh-cn.xlf"

$zhNameA2 = "61df56d7-e540-4367-8ba3-57a9d998e113.md"
$zhNameD2 = "61df56d7-e540-4367-8ba3-57a9d998e113.2c65cb27da0e9d1c9c46324e0027fa7f7b8ec244.zh-cn.xlf"
$zhNameA3 = "b3c9e954-dfb6-437b-b587-b10869f38c87.md"
$zhNameD3 = "b3c9e954-dfb6-437b-b587-b10869f38c87.264bcc2d9140a76805ec1f79e3eaf0d34992c5a9.zh-cn.xlf"

# New handback datetime for the zh-cn language pass.
$wsZhCn.Range("H2").Value2 = "2016-03-13 11:02:41"
$wsZhCn.Range("H3").Value2 = "2016-03-13 11:02:41"

# Rebuild every hyperlink on the sheet, in left-to-right / top-to-bottom
# order, so the newly-introduced F/G links land between the existing D and
# A(next row) links exactly like a freshly generated report would.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhUrlA2, "", "", $zhNameA2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $zhUrlA2, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhUrlD2, "", "", $zhNameD2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhUrlA2, "", "", $zhNameA2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhUrlD2, "", "", $zhNameD2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $zhUrlA3, "", "", $zhNameA3)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $zhUrlA3, "", "", ".md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $zhUrlD3, "", "", $zhNameD3)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhUrlA3, "", "", $zhNameA3)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhUrlD3, "", "", $zhNameD3)

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape, its own URLs/datetime.
# ---------------------------------------------------------------------------
$deUrlA2 = "https://github.com/OpenLocalizationTest/oltest/blob/6d562aac7f97325115fbb32b09ff899fb3811720/e2e/61df56d7-e540-4367-8ba3-57a9d998e113.md"
$deUrlD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce16dc16e0c0978f3f8b95bd04a964bb7ced6b5e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/61df56d7-e540-4367-8ba3-57a9d998e113.2c65cb27da0e9d1c9c46324e0027fa7f7b8ec244.de-de.xlf"
$deUrlA3 = "https://github.com/OpenLocalizationTest/oltest/blob/6d562aac7f97325115fbb32b09ff899fb3811720/e2e/b3c9e954-dfb6-437b-b587-b10869f38c87.md"
$deUrlD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ce16dc16e0c0978f3f8b95bd04a964bb7ced6b5e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b3c9e954-dfb6-437b-b587-b10869f38c87.264bcc2d9140a76805ec1f79e3eaf0d34992c5a9.de-de.xlf"

$deNameA2 = "61df56d7-e540-4367-8ba3-57a9d998e113.md"
$deNameD2 = "61df56d7-e540-4367-8ba3-57a9d998e113.2c65cb27da0e9d1c9c46324e0027fa7f7b8ec244.de-de.xlf"
$deNameA3 = "b3c9e954-dfb6-437b-b587-b10869f38c87.md"
$deNameD3 = "b3c9e954-dfb6-437b-b587-b10869f38c87.264bcc2d9140a76805ec1f79e3eaf0d34992c5a9.de-de.xlf"

$wsDeDe.Range("H2").Value2 = "2016-03-13 11:02:48"
$wsDeDe.Range("H3").Value2 = "2016-03-13 11:02:48"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deUrlA2, "", "", $deNameA2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $deUrlA2, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deUrlD2, "", "", $deNameD2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deUrlA2, "", "", $deNameA2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deUrlD2, "", "", $deNameD2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $deUrlA3, "", "", $deNameA3)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $deUrlA3, "", "", ".md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $deUrlD3, "", "", $deNameD3)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deUrlA3, "", "", $deNameA3)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deUrlD3, "", "", $deNameD3)
